$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-6 from 2023-10-08 (45207) to 2023-10-09 (45208)
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 3).Value = 45208
}
